# Generate Report for Handoff
# Refresh the localization-status report: 3 old files are replaced by 4 new
# files (new GUIDs, new hashes), all now ".md" sources, all "Ready for
# handoff" / "Include", with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$baseSha = "ae075f8e7bdc0caf0b2d172efb22789c08e07373"
$zhSha   = "312500e4403b4a1501f96cfbfc2cd9a542564142"
$deSha   = "d8d35b004ebd86c446d7a1e050176422e0dfa2c4"

# The four source files now tracked by the report.
$files = @(
    @{ name = "5c6ad8be-cea4-4e1e-b77a-f1f5fc22ebed.md"; hash = "706137c8a29d1c102b8d381e753225a659fed7e5" },
    @{ name = "6b013868-cb30-4de4-974c-07e834b919e2.md"; hash = "88e565685c5b4b5087a3f433c6b54c864e1e9e4a" },
    @{ name = "7e4bcb88-d4a4-48d9-971d-2c44acdfd2dd.md"; hash = "149999e07c750f7332f7dea1b9dac9f1951f7df4" },
    @{ name = "d3a5b5e8-7cd3-4563-b175-1cf3fe2d71bc.md"; hash = "74195a73e1909a21cef359f14241e009d94b6701" }
)

$overviewDate  = "2016-03-19 04:05:53"
$handoffDate   = "2016-03-19 04:05:44"
$epoch         = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview" : File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Drop the old hyperlinks before rewriting the rows underneath them.
$wsOv.Range("A1:D20").Hyperlinks.Delete()

$row = 2
foreach ($f in $files) {
    $wsOv.Range("A$row").Value = $f.name
    $wsOv.Range("B$row").Value = "Ready for handoff"
    $wsOv.Range("C$row").Value = "Ready for handoff"
    $wsOv.Range("D$row").Value = $overviewDate
    $wsOv.Range("D$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $url = "https://github.com/OpenLocalizationTest/oltest/blob/$baseSha/e2e/$($f.name)"
    $wsOv.Hyperlinks.Add($wsOv.Range("A$row"), $url, "", "", $f.name)

    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de" : per-language handoff detail
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $lang, $repoSha, $langDate) {
    $ws.Range("A1:L20").Hyperlinks.Delete()

    $row = 2
    foreach ($f in $files) {
        $xlf = "$($f.name).$($f.hash).$lang.xlf"

        $ws.Range("A$row").Value = $f.name
        $ws.Range("B$row").Value = ".md"
        $ws.Range("C$row").Value = "Ready for handoff"
        $ws.Range("D$row").Value = $xlf
        $ws.Range("E$row").Value = $langDate
        $ws.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("H$row").Value = $epoch
        $ws.Range("H$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("I$row").Value = ""
        $ws.Range("J$row").Value = "Include"

        $srcUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$baseSha/e2e/$($f.name)"
        $ws.Hyperlinks.Add($ws.Range("A$row"), $srcUrl, "", "", $f.name)

        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$repoSha/ol-handoff/OpenLocalizationTestOrg/oltest-$($lang -replace '-','')-fly/yuwzho/ht/$xlf"
        $ws.Hyperlinks.Add($ws.Range("D$row"), $xlfUrl, "", "", $xlf)

        $row = $row + 1
    }
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZh "zh-cn" $zhSha $handoffDate

$wsDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDe "de-de" $deSha $overviewDate
